# Apply the price/coin-listing updates for the 16-12-2022 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.16"
$ws.Range("D2").Style = "Normal"

$ws.Range("D4").Value = "'6.202"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.06162"
$ws.Range("D5").Style = "Normal"

$ws.Range("D7").Value = "'3.466"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'1.346"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.7992"
$ws.Range("D9").Style = "Normal"

$ws.Range("B10").Value = "WazirX"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1577"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.08087"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.03509"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.03092"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09311"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "MCDex"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'3.841"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.001693"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.04800"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "One"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'0.0006139"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.006195"
$ws.Range("D19").Style = "Normal"

$ws.Range("D21").Value = "'0.004073"
$ws.Range("D21").Style = "Normal"

$ws.Range("D23").Value = "'3.693"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "'2.214"
$ws.Range("D24").Style = "Normal"

$ws.Range("D41").Value = "'0.007141"
$ws.Range("D41").Style = "Normal"

$ws.Range("D43").Value = "'0.003130"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.01001"
$ws.Range("D44").Style = "Normal"

$ws.Range("D46").Value = "'0.00005940"
$ws.Range("D46").Style = "Normal"

$ws.Range("D48").Value = "'0.6999"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "'0.1521"
$ws.Range("D49").Style = "Normal"
